$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 (Oppo Service Center) updates, mirroring a Streamlit app data refresh:
#  - Start Date (B9) moves forward one year
#  - Lease End Date (C9) moves forward one year
#  - 3-Month Reminder (D9) cleared since lease is no longer ending soon
#  - Projected Income (F9) now matches Actual Income (G9)

$ws.Range("B9").Value = 45822
$ws.Range("C9").Value = 46187
$ws.Range("D9").Value = ""
$ws.Range("F9").Value = 81000000
